$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 199.25
$ws.Range("I9").Value = 199.25
$ws.Range("K9").Value = 199.25
$ws.Range("M9").Value = -30.25
$ws.Range("H70").Value = 3482
$ws.Range("J70").Value = 3784.3635
$ws.Range("L70").Value = 11353.0905
$ws.Range("N70").Value = -11893.0905
$ws.Range("H73").Value = 3482
$ws.Range("J73").Value = 3784.3635
$ws.Range("L73").Value = 11353.0905
$ws.Range("N73").Value = -13225.0905
$ws.Range("H87").Value = 58500
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 58500
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 58500
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -60996
$ws.Range("H90").Value = 58500
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 58500
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 175500
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -187980
$ws.Range("H115").Value = 700
$ws.Range("I115").Value = 700
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 2100
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -533
$ws.Range("N115").ClearContents()
$ws.Range("H137").Value = 2513.8572
$ws.Range("I137").Value = 2051
$ws.Range("K137").Value = 6153
$ws.Range("M137").Value = -3603
$ws.Range("H138").Value = 2693.3062
$ws.Range("I138").Value = 2038.1666
$ws.Range("J138").Value = 3073.7097
$ws.Range("K138").Value = 6114.4998
$ws.Range("L138").Value = 9221.1291
$ws.Range("M138").Value = -974.4997999999996
$ws.Range("N138").Value = -19501.1291

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 3408.111
$ws.Range("I5").Value = 4342.5713
$ws.Range("K5").Value = 4342.5713
$ws.Range("M5").Value = -4230.5713
$ws.Range("H45").Value = 49162.09
$ws.Range("I45").Value = 66833.74000000001
$ws.Range("K45").Value = 66833.74000000001
$ws.Range("M45").Value = -66456.74000000001
$ws.Range("H132").Value = 4268.74
$ws.Range("I132").Value = 3429.361
$ws.Range("K132").Value = 10288.083
$ws.Range("M132").Value = -7758.082999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 3408.111
$ws.Range("I4").Value = 4342.5713
$ws.Range("K4").Value = 4342.5713
$ws.Range("M4").Value = -4227.5713
$ws.Range("H76").Value = 2901.6667
$ws.Range("J76").Value = 2901.6667
$ws.Range("L76").Value = 2901.6667
$ws.Range("N76").Value = -3531.6667
$ws.Range("H79").Value = 2901.6667
$ws.Range("J79").Value = 2901.6667
$ws.Range("L79").Value = 2901.6667
$ws.Range("N79").Value = -5085.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 30296.334
$ws.Range("I51").Value = 30289
$ws.Range("K51").Value = 30289
$ws.Range("M51").Value = -29553
$ws.Range("H61").Value = 30296.334
$ws.Range("I61").Value = 30289
$ws.Range("K61").Value = 30289
$ws.Range("M61").Value = -29941
$ws.Range("H74").Value = 41920
$ws.Range("J74").Value = 41920
$ws.Range("L74").Value = 41920
$ws.Range("N74").Value = -43668
$ws.Range("H77").Value = 41920
$ws.Range("J77").Value = 41920
$ws.Range("L77").Value = 125760
$ws.Range("N77").Value = -134496
$ws.Range("H94").Value = 4864.8
$ws.Range("J94").Value = 4864.8
$ws.Range("L94").Value = 4864.8
$ws.Range("N94").Value = -5766.8
$ws.Range("H97").Value = 25611.5
$ws.Range("J97").Value = 25611.5
$ws.Range("L97").Value = 25611.5
$ws.Range("N97").Value = -27593.5
$ws.Range("H99").Value = 8597.6
$ws.Range("I99").Value = 6995.6665
$ws.Range("J99").Value = 9284.143
$ws.Range("K99").Value = 6995.6665
$ws.Range("L99").Value = 9284.143
$ws.Range("M99").Value = -5497.6665
$ws.Range("N99").Value = -12280.143
$ws.Range("H126").Value = 8597.6
$ws.Range("I126").Value = 6995.6665
$ws.Range("J126").Value = 9284.143
$ws.Range("K126").Value = 20986.9995
$ws.Range("L126").Value = 27852.429
$ws.Range("M126").Value = -18516.9995
$ws.Range("N126").Value = -32792.429
$ws.Range("H132").Value = 2572.7778
$ws.Range("I132").Value = 2306.5715
$ws.Range("J132").Value = 3504.5
$ws.Range("K132").Value = 6919.7145
$ws.Range("L132").Value = 10513.5
$ws.Range("M132").Value = -4389.7145
$ws.Range("N132").Value = -15573.5
$ws.Range("H141").Value = 37900
$ws.Range("J141").Value = 37900
$ws.Range("L141").Value = 37900
$ws.Range("N141").Value = -48260

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 4596.625
$ws.Range("I2").Value = 491
$ws.Range("J2").Value = 25124.75
$ws.Range("K2").Value = 2946
$ws.Range("L2").Value = 150748.5
$ws.Range("M2").Value = -2833
$ws.Range("N2").Value = -150974.5
$ws.Range("H28").Value = 2910
$ws.Range("I28").Value = 3365
$ws.Range("J28").Value = 2000
$ws.Range("K28").Value = 10095
$ws.Range("L28").Value = 6000
$ws.Range("M28").Value = -9863
$ws.Range("N28").Value = -6464
$ws.Range("H128").Value = 747892.25
$ws.Range("I128").Value = 747892.25
$ws.Range("K128").Value = 2243676.75
$ws.Range("M128").Value = -2238696.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 30852.428
$ws.Range("I46").Value = 22649.5
$ws.Range("J46").Value = 34133.6
$ws.Range("K46").Value = 22649.5
$ws.Range("L46").Value = 34133.6
$ws.Range("M46").Value = -22493.5
$ws.Range("N46").Value = -34445.6
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H122").Value = 2105.8696
$ws.Range("I122").Value = 2208.0715
$ws.Range("J122").Value = 1946.8889
$ws.Range("K122").Value = 6624.2145
$ws.Range("L122").Value = 5840.6667
$ws.Range("M122").Value = -4174.2145
$ws.Range("N122").Value = -10740.6667
$ws.Range("H132").Value = 7431.227
$ws.Range("I132").Value = 7332.778
$ws.Range("J132").Value = 7874.25
$ws.Range("K132").Value = 21998.334
$ws.Range("L132").Value = 23622.75
$ws.Range("M132").Value = -19468.334
$ws.Range("N132").Value = -28682.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8921.6
$ws.Range("I7").Value = 8921.6
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 8921.6
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -8809.6
$ws.Range("N7").ClearContents()
$ws.Range("H20").Value = 20619.375
$ws.Range("J20").Value = 20619.375
$ws.Range("L20").Value = 20619.375
$ws.Range("N20").Value = -21071.375
$ws.Range("H22").Value = 798.75
$ws.Range("I22").Value = 945
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 945
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -650
$ws.Range("N22").Value = -1340
$ws.Range("H27").Value = 798.75
$ws.Range("I27").Value = 945
$ws.Range("J27").Value = 750
$ws.Range("K27").Value = 945
$ws.Range("L27").Value = 750
$ws.Range("M27").Value = -838
$ws.Range("N27").Value = -964
$ws.Range("H75").Value = 16891
$ws.Range("J75").Value = 16891
$ws.Range("L75").Value = 16891
$ws.Range("N75").Value = -18763
$ws.Range("H78").Value = 16891
$ws.Range("J78").Value = 16891
$ws.Range("L78").Value = 50673
$ws.Range("N78").Value = -60033
$ws.Range("H122").Value = 4211.2354
$ws.Range("J122").Value = 7099.2
$ws.Range("L122").Value = 21297.6
$ws.Range("N122").Value = -26197.6
$ws.Range("H126").Value = 8921.6
$ws.Range("I126").Value = 8921.6
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 26764.8
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -24294.8
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 7193.645
$ws.Range("I132").Value = 7378.885
$ws.Range("J132").Value = 6230.4
$ws.Range("K132").Value = 22136.655
$ws.Range("L132").Value = 18691.2
$ws.Range("M132").Value = -19606.655
$ws.Range("N132").Value = -23751.2
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -75060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 38899.8
$ws.Range("J70").Value = 38899.8
$ws.Range("L70").Value = 38899.8
$ws.Range("N70").Value = -39529.8
$ws.Range("H73").Value = 38899.8
$ws.Range("J73").Value = 38899.8
$ws.Range("L73").Value = 38899.8
$ws.Range("N73").Value = -41083.8
$ws.Range("H81").Value = 4669.25
$ws.Range("I81").Value = 2750.5
$ws.Range("J81").Value = 6588
$ws.Range("K81").Value = 5501
$ws.Range("L81").Value = 13176
$ws.Range("M81").Value = -4440
$ws.Range("N81").Value = -15298
$ws.Range("H84").Value = 4669.25
$ws.Range("I84").Value = 2750.5
$ws.Range("J84").Value = 6588
$ws.Range("K84").Value = 27505
$ws.Range("L84").Value = 65880
$ws.Range("M84").Value = -22201
$ws.Range("N84").Value = -76488
$ws.Range("H126").Value = 5796.759
$ws.Range("I126").Value = 5836.28
$ws.Range("K126").Value = 17508.84
$ws.Range("M126").Value = -15038.84
$ws.Range("H132").Value = 7812.5264
$ws.Range("I132").Value = 6562.6
$ws.Range("J132").Value = 12499.75
$ws.Range("K132").Value = 19687.8
$ws.Range("L132").Value = 37499.25
$ws.Range("M132").Value = -17157.8
$ws.Range("N132").Value = -42559.25
